# Optimize navigation menu accessibility
# Append a new data row (row 23) to each of the four worksheets in the
# SAG1 database workbook, mirroring the structure of the existing rows.

$wb = $excel.ActiveWorkbook

# Sheet name -> new row values
# Columns: A time, B total-len(hex), C id(hex), D actual-len(hex), E checksum(hex),
#          F total-len dec, G id dec, H actual-len dec, I checksum dec
$newRows = @{
    "ROW35-FE-LIFTER"  = @{
        A = 45735.77978195602
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x82"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 386
        I = 13
    }
    "ROW35-MID-LIFTER" = @{
        A = 45735.62962542824
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x82"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 386
        I = 14
    }
    "ROW02-FE-LIFTER"  = @{
        A = 45735.77790322917
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x82"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 386
        I = 3
    }
    "ROW02-MID-LIFTER" = @{
        A = 45735.83884072916
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x82"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 386
        I = 3
    }
}

foreach ($sheetName in $newRows.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = $newRows[$sheetName]

    $targetRow = 23

    $ws.Cells.Item($targetRow, 1).Value = $row.A
    $ws.Cells.Item($targetRow, 1).NumberFormat = $ws.Cells.Item($targetRow - 1, 1).NumberFormat

    $ws.Cells.Item($targetRow, 2).Value = $row.B
    $ws.Cells.Item($targetRow, 3).Value = $row.C
    $ws.Cells.Item($targetRow, 4).Value = $row.D
    $ws.Cells.Item($targetRow, 5).Value = $row.E
    $ws.Cells.Item($targetRow, 6).Value = $row.F
    $ws.Cells.Item($targetRow, 7).Value = $row.G
    $ws.Cells.Item($targetRow, 8).Value = $row.H
    $ws.Cells.Item($targetRow, 9).Value = $row.I
}
